$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped from the commit diff (coin prices/volumes refreshed,
# plus three coins re-ranked: XRP<->BNB swap rows 5/6, and a 3-way rotation
# THORChain/ARBITRUM/VeChain across rows 40/41/42).
# Values that look numeric are apostrophe-prefixed so Excel stores them as
# literal text (matching the original inlineStr/text cells) instead of
# auto-converting them to numbers.
$updates = @(
    @{Cell='D2'; Value='36.844.18'}
    @{Cell='E2'; Value='  -0.84%  '}
    @{Cell='D3'; Value='2.094.97'}
    @{Cell='E3'; Value='  +2.17%  '}
    @{Cell='E4'; Value='  -0.08%  '}
    @{Cell='B5'; Value='XRP'}
    @{Cell='C5'; Value='https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'}
    @{Cell='D5'; Value='''0.703'}
    @{Cell='E5'; Value='  +5.73%  '}
    @{Cell='B6'; Value='BNB'}
    @{Cell='C6'; Value='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'}
    @{Cell='D6'; Value='''245.74'}
    @{Cell='E6'; Value='  -0.98%  '}
    @{Cell='E7'; Value='  -0.01%  '}
    @{Cell='D8'; Value='''54.06'}
    @{Cell='E8'; Value='  -4.87%  '}
    @{Cell='D9'; Value='''59.21'}
    @{Cell='E9'; Value='  -1.42%  '}
    @{Cell='D10'; Value='''0.369'}
    @{Cell='E10'; Value='  -3.80%  '}
    @{Cell='D11'; Value='''0.0770'}
    @{Cell='E11'; Value='  -1.71%  '}
    @{Cell='E12'; Value='  +1.37%  '}
    @{Cell='D13'; Value='''0.922'}
    @{Cell='E13'; Value='  +0.71%  '}
    @{Cell='D14'; Value='''14.97'}
    @{Cell='E14'; Value='  -7.63%  '}
    @{Cell='D15'; Value='2.399.54'}
    @{Cell='E15'; Value='  +2.14%  '}
    @{Cell='D16'; Value='''5.50'}
    @{Cell='E16'; Value='  -4.28%  '}
    @{Cell='D17'; Value='2.117.48'}
    @{Cell='E17'; Value='  +3.30%  '}
    @{Cell='D18'; Value='36.818.05'}
    @{Cell='E18'; Value='  -0.90%  '}
    @{Cell='D19'; Value='''17.23'}
    @{Cell='E19'; Value='  -8.10%  '}
    @{Cell='D20'; Value='''72.86'}
    @{Cell='E20'; Value='  -2.25%  '}
    @{Cell='D21'; Value='0.0₃0885'}
    @{Cell='E21'; Value='  -1.47%  '}
    @{Cell='E22'; Value='  +0.19%  '}
    @{Cell='D23'; Value='''239.65'}
    @{Cell='E23'; Value='  +1.14%  '}
    @{Cell='E24'; Value='  +0.09%  '}
    @{Cell='D25'; Value='''2.39'}
    @{Cell='E25'; Value='  -3.49%  '}
    @{Cell='D26'; Value='''9.65'}
    @{Cell='E26'; Value='  +0.01%  '}
    @{Cell='E27'; Value='  -0.97%  '}
    @{Cell='D28'; Value='''167.23'}
    @{Cell='E28'; Value='  -1.64%  '}
    @{Cell='D29'; Value='''20.94'}
    @{Cell='E29'; Value='  +3.64%  '}
    @{Cell='D30'; Value='''0.127'}
    @{Cell='E30'; Value='  +1.41%  '}
    @{Cell='D31'; Value='''5.25'}
    @{Cell='E31'; Value='  +4.35%  '}
    @{Cell='D32'; Value='''1.17'}
    @{Cell='E32'; Value='  +0.10%  '}
    @{Cell='D33'; Value='''4.74'}
    @{Cell='E33'; Value='  +4.58%  '}
    @{Cell='D34'; Value='''0.0609'}
    @{Cell='E34'; Value='  -2.15%  '}
    @{Cell='E35'; Value='  +8.06%  '}
    @{Cell='E36'; Value='  -0.04%  '}
    @{Cell='E37'; Value='  +3.65%  '}
    @{Cell='E38'; Value='  -6.52%  '}
    @{Cell='E39'; Value='  -4.84%  '}
    @{Cell='B40'; Value='THORChain'}
    @{Cell='C40'; Value='https://coinranking.com/coin/ybmU-kKU+thorchain-rune'}
    @{Cell='D40'; Value='''4.92'}
    @{Cell='E40'; Value='  -7.32%  '}
    @{Cell='B41'; Value='ARBITRUM'}
    @{Cell='C41'; Value='https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'}
    @{Cell='D41'; Value='''1.16'}
    @{Cell='E41'; Value='  +0.90%  '}
    @{Cell='B42'; Value='VeChain'}
    @{Cell='C42'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'}
    @{Cell='D42'; Value='''0.0221'}
    @{Cell='E42'; Value='  -1.35%  '}
    @{Cell='D43'; Value='''0.0964'}
    @{Cell='E43'; Value='  -1.52%  '}
    @{Cell='D44'; Value='''96.55'}
    @{Cell='D45'; Value='''2.87'}
    @{Cell='E45'; Value='  -7.26%  '}
    @{Cell='D46'; Value='''7.87'}
    @{Cell='E46'; Value='  +15.43%  '}
    @{Cell='D47'; Value='1.410.22'}
    @{Cell='D48'; Value='''16.09'}
    @{Cell='E48'; Value='  -8.48%  '}
    @{Cell='D49'; Value='''2.43'}
    @{Cell='E49'; Value='  -0.42%  '}
    @{Cell='D50'; Value='''2.90'}
    @{Cell='E50'; Value='  +1.51%  '}
    @{Cell='D51'; Value='2.287.29'}
    @{Cell='E51'; Value='  +2.24%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
